$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 145 and 146: fill in resultado/profit ---
$ws.Cells.Item(145, 7).Value = "Fallo"
$ws.Cells.Item(145, 8).Value = -1

$ws.Cells.Item(146, 7).Value = "Fallo"
$ws.Cells.Item(146, 8).Value = -1

# --- Append new rows 147-153 ---
# Helper approach: for date-like text in column B we must force Text
# formatting BEFORE assigning the value, otherwise Excel auto-converts
# "yyyy-mm-dd" strings into date serial numbers. Resetting the style
# afterwards keeps the cell a plain (unstyled) text cell.

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 147
$ws.Cells.Item(147, 1).Value = 14870662
Set-TextCell 147 2 "2025-10-13"
$ws.Cells.Item(147, 3).Value = "Roberto Bautista Agut"
$ws.Cells.Item(147, 4).Value = "Sebastián Báez"
$ws.Cells.Item(147, 5).Value = "Gana Sebastián Báez"
$ws.Cells.Item(147, 6).Value = 3

# Row 148
$ws.Cells.Item(148, 1).Value = 14871770
Set-TextCell 148 2 "2025-10-13"
$ws.Cells.Item(148, 3).Value = "Tallon Griekspoor"
$ws.Cells.Item(148, 4).Value = "Jacob Fearnley"
$ws.Cells.Item(148, 5).Value = "Gana Tallon Griekspoor"
$ws.Cells.Item(148, 6).Value = 1.73

# Row 149
$ws.Cells.Item(149, 1).Value = 14871728
Set-TextCell 149 2 "2025-10-13"
$ws.Cells.Item(149, 3).Value = "Shuai Zhang"
$ws.Cells.Item(149, 4).Value = "Veronika Kudermetova"
$ws.Cells.Item(149, 5).Value = "Gana Veronika Kudermetova"
$ws.Cells.Item(149, 6).Value = 1.57

# Row 150
$ws.Cells.Item(150, 1).Value = 14871732
Set-TextCell 150 2 "2025-10-13"
$ws.Cells.Item(150, 3).Value = "Yue Yuan"
$ws.Cells.Item(150, 4).Value = "Xinyu Wang"
$ws.Cells.Item(150, 5).Value = "Gana Xinyu Wang"
$ws.Cells.Item(150, 6).Value = 1.73

# Row 151
$ws.Cells.Item(151, 1).Value = 14870698
Set-TextCell 151 2 "2025-10-12"
$ws.Cells.Item(151, 3).Value = "Marc-Andrea Huesler"
$ws.Cells.Item(151, 4).Value = "Bernard Tomic"
$ws.Cells.Item(151, 5).Value = "Gana Bernard Tomic"
$ws.Cells.Item(151, 6).Value = 2.38

# Row 152
$ws.Cells.Item(152, 1).Value = 14870026
Set-TextCell 152 2 "2025-10-13"
$ws.Cells.Item(152, 3).Value = "Naoya Honda"
$ws.Cells.Item(152, 4).Value = "Frederico Ferreira Silva"
$ws.Cells.Item(152, 5).Value = "Gana Naoya Honda"
$ws.Cells.Item(152, 6).Value = 3.5

# Row 153
$ws.Cells.Item(153, 1).Value = 14870031
Set-TextCell 153 2 "2025-10-13"
$ws.Cells.Item(153, 3).Value = "Hyeon Chung"
$ws.Cells.Item(153, 4).Value = "Te Rigele"
$ws.Cells.Item(153, 5).Value = "Gana Te Rigele"
$ws.Cells.Item(153, 6).Value = 2.38
